$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 664.8570999999999
$ws.Range("I32").Value = 400
$ws.Range("J32").Value = 770.8
$ws.Range("K32").Value = 400
$ws.Range("L32").Value = 770.8
$ws.Range("M32").Value = -74
$ws.Range("N32").Value = -1422.8
$ws.Range("H43").Value = 1800
$ws.Range("I43").Value = 1500
$ws.Range("J43").Value = 1860
$ws.Range("K43").Value = 1500
$ws.Range("L43").Value = 1860
$ws.Range("M43").Value = -1431
$ws.Range("N43").Value = -1998
$ws.Range("H64").Value = 3072.6667
$ws.Range("I64").Value = 2864
$ws.Range("K64").Value = 2864
$ws.Range("M64").Value = -2616
$ws.Range("H67").Value = 3072.6667
$ws.Range("I67").Value = 2864
$ws.Range("K67").Value = 2864
$ws.Range("M67").Value = -2006
$ws.Range("H94").Value = 2997.1428
$ws.Range("H111").Value = 4004.1428
$ws.Range("I111").Value = 4205.8
$ws.Range("J111").Value = 3500
$ws.Range("K111").Value = 12617.4
$ws.Range("L111").Value = 10500
$ws.Range("M111").Value = -9550.400000000001
$ws.Range("N111").Value = -16634
$ws.Range("H130").Value = 39780
$ws.Range("J130").Value = 39780
$ws.Range("L130").Value = 39780
$ws.Range("N130").Value = -49820
$ws.Range("H138").Value = 2565.2715
$ws.Range("J138").Value = 3057.9805
$ws.Range("L138").Value = 9173.941500000001
$ws.Range("N138").Value = -19453.9415
$ws.Range("H141").Value = 2498.3333
$ws.Range("I141").Value = 1421.1538
$ws.Range("J141").Value = 9500
$ws.Range("K141").Value = 4263.4614
$ws.Range("L141").Value = 28500
$ws.Range("M141").Value = 916.5385999999999
$ws.Range("N141").Value = -38860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1570.2106
$ws.Range("I61").Value = 1501.8889
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 1501.8889
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -1289.8889
$ws.Range("N61").Value = -3224
$ws.Range("H136").Value = 1570.2106
$ws.Range("I136").Value = 1501.8889
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 4505.6667
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -1955.6667
$ws.Range("N136").Value = -13500

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 26181.715
$ws.Range("J95").Value = 26181.715
$ws.Range("L95").Value = 26181.715
$ws.Range("N95").Value = -31673.715
$ws.Range("H107").Value = 1698.8235
$ws.Range("I107").Value = 1411.8182
$ws.Range("J107").Value = 2225
$ws.Range("K107").Value = 1411.8182
$ws.Range("L107").Value = 2225
$ws.Range("M107").Value = 508.1818000000001
$ws.Range("N107").Value = -6065

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 4736.5835
$ws.Range("J44").Value = 5084.1
$ws.Range("L44").Value = 15252.3
$ws.Range("N44").Value = -16048.3
$ws.Range("H51").Value = 4133.684
$ws.Range("J51").Value = 4352.222
$ws.Range("L51").Value = 13056.666
$ws.Range("N51").Value = -13976.666
$ws.Range("H122").Value = 735.9
$ws.Range("I122").Value = 330
$ws.Range("K122").Value = 2970
$ws.Range("M122").Value = -520

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 28000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 28000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 28000
$ws.Range("N69").Value = -29498
$ws.Range("H72").Value = 28000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 28000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 84000
$ws.Range("N72").Value = -91488
$ws.Range("H122").Value = 2507.8333
$ws.Range("I122").Value = 2809.4
$ws.Range("K122").Value = 8428.200000000001
$ws.Range("M122").Value = -5978.200000000001
$ws.Range("H123").Value = 17546.441
$ws.Range("J123").Value = 17546.441
$ws.Range("L123").Value = 17546.441
$ws.Range("N123").Value = -22446.441
$ws.Range("H126").Value = 2643.3333
$ws.Range("I126").Value = 1741.6666
$ws.Range("J126").Value = 3725.3333
$ws.Range("K126").Value = 5224.9998
$ws.Range("L126").Value = 11175.9999
$ws.Range("M126").Value = -2754.9998
$ws.Range("N126").Value = -16115.9999
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1578.9445
$ws.Range("I7").Value = 1072.8182
$ws.Range("J7").Value = 1801.64
$ws.Range("K7").Value = 1072.8182
$ws.Range("L7").Value = 1801.64
$ws.Range("M7").Value = -960.8181999999999
$ws.Range("N7").Value = -2025.64
$ws.Range("H46").Value = 977.375
$ws.Range("I46").Value = 575
$ws.Range("J46").Value = 1111.5
$ws.Range("K46").Value = 575
$ws.Range("L46").Value = 1111.5
$ws.Range("M46").Value = -387
$ws.Range("N46").Value = -1487.5
$ws.Range("H126").Value = 1578.9445
$ws.Range("I126").Value = 1072.8182
$ws.Range("J126").Value = 1801.64
$ws.Range("K126").Value = 3218.4546
$ws.Range("L126").Value = 5404.92
$ws.Range("M126").Value = -748.4546
$ws.Range("N126").Value = -10344.92
$ws.Range("H135").Value = 55142.668
$ws.Range("J135").Value = 55142.668
$ws.Range("L135").Value = 55142.668
$ws.Range("N135").Value = -65282.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1305
$ws.Range("I23").Value = 1305
$ws.Range("K23").Value = 1305
$ws.Range("M23").Value = -1076
$ws.Range("H43").Value = 6924.75
$ws.Range("I43").Value = 4233
$ws.Range("K43").Value = 4233
$ws.Range("M43").Value = -4084
$ws.Range("H52").Value = 16666.666
$ws.Range("J52").Value = 21000
$ws.Range("L52").Value = 21000
$ws.Range("N52").Value = -21452
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("H126").Value = 2691.8572
$ws.Range("I126").Value = 2807.1667
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 8421.500100000001
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -5951.500100000001
$ws.Range("N126").Value = -10940
$ws.Range("M58").ClearContents()
